# Commit: "#5: property boat&car done"
# The 汽車 (car) sheet's row 1 was actually a stray duplicate of the data
# row instead of a proper header, and the data row (row 2) was missing the
# trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that every other
# property sheet in this workbook already carries. This rebuilds row 1 as
# the real header and fills out row 2 with the missing columns, matching
# the 14-column (A:N) layout used by the other sheets (e.g. 土地/建物).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper header labels (B1:N1) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Make sure the new header cells carry the same (bold/centered/bordered)
# header style already used on this row.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2: fix the car name/register_date values, fill in the rest ---
$ws.Range("B2").Value = "豐田國瑞"
$ws.Range("C2").Value = 1998
$ws.Range("D2").Value = "王廷升"
$ws.Range("E2").Value = "95年07月21日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 700000
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# "date" is the source-file timestamp string "2012-04-30" (matches the
# other sheets' J/date column, stored as literal text there too) -- force
# text so Excel doesn't reinterpret the dashed digits as a real date.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("K2").Value = "王廷升"
$ws.Range("L2").Value = 1727
$ws.Range("M2").Value = "tmpc32d1"
$ws.Range("N2").Value = 33

# Carry the plain data-row style onto the newly added cells too.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
